$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 12).Value = 576.1111
$ws.Cells.Item(4, 9).Value = 45.714287
$ws.Cells.Item(4, 10).Value = 576.1111
$ws.Cells.Item(4, 13).Value = 68.285713
$ws.Cells.Item(4, 14).Value = -804.1111
$ws.Cells.Item(4, 11).Value = 45.714287
$ws.Cells.Item(4, 8).Value = 344.0625
$ws.Cells.Item(8, 11).Value = 367.5
$ws.Cells.Item(8, 8).Value = 122.5
$ws.Cells.Item(8, 13).Value = -228.5
$ws.Cells.Item(8, 9).Value = 122.5
$ws.Cells.Item(16, 8).Value = 27722
$ws.Cells.Item(16, 12).Value = 7444
$ws.Cells.Item(16, 10).Value = 7444
$ws.Cells.Item(16, 14).Value = -7904
$ws.Cells.Item(51, 8).Value = 7337.5
$ws.Cells.Item(51, 12).Value = 7337.5
$ws.Cells.Item(51, 10).Value = 7337.5
$ws.Cells.Item(51, 14).Value = -8305.5
$ws.Cells.Item(86, 11).Value = 4240.5
$ws.Cells.Item(86, 8).Value = 4874.2915
$ws.Cells.Item(86, 13).Value = -3117.5
$ws.Cells.Item(86, 9).Value = 4240.5
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 14).ClearContents()
$ws.Cells.Item(89, 11).Value = 21202.5
$ws.Cells.Item(89, 8).Value = 4874.2915
$ws.Cells.Item(89, 9).Value = 4240.5
$ws.Cells.Item(89, 13).Value = -15586.5
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 14).ClearContents()
$ws.Cells.Item(107, 9).Value = 27862088
$ws.Cells.Item(107, 13).Value = -27860168
$ws.Cells.Item(107, 11).Value = 27862088
$ws.Cells.Item(107, 8).Value = 19668036
$ws.Cells.Item(112, 11).Value = 6984.999899999999
$ws.Cells.Item(112, 12).Value = 24573.375
$ws.Cells.Item(112, 8).Value = 7265.421
$ws.Cells.Item(112, 13).Value = -5876.999899999999
$ws.Cells.Item(112, 9).Value = 2328.3333
$ws.Cells.Item(112, 10).Value = 8191.125
$ws.Cells.Item(112, 14).Value = -26789.375
$ws.Cells.Item(113, 8).Value = 7906.154
$ws.Cells.Item(113, 12).Value = 7616.364
$ws.Cells.Item(113, 10).Value = 7616.364
$ws.Cells.Item(113, 14).Value = -14124.364
$ws.Cells.Item(116, 8).Value = 5966.125
$ws.Cells.Item(116, 12).Value = 6104.857
$ws.Cells.Item(116, 10).Value = 6104.857
$ws.Cells.Item(116, 14).Value = -12988.857
$ws.Cells.Item(121, 12).Value = 7703.000100000001
$ws.Cells.Item(121, 8).Value = 2544.96
$ws.Cells.Item(121, 10).Value = 2567.6667
$ws.Cells.Item(121, 14).Value = -11197.0001
$ws.Cells.Item(125, 8).Value = 5466959
$ws.Cells.Item(125, 13).Value = -11928.75
$ws.Cells.Item(125, 12).Value = 52654437
$ws.Cells.Item(125, 9).Value = 1598.75
$ws.Cells.Item(125, 10).Value = 5850493
$ws.Cells.Item(125, 14).Value = -52659357
$ws.Cells.Item(125, 11).Value = 14388.75
$ws.Cells.Item(132, 8).Value = 10872383
$ws.Cells.Item(132, 12).Value = 11430.9999
$ws.Cells.Item(132, 9).Value = 11630656
$ws.Cells.Item(132, 10).Value = 3810.3333
$ws.Cells.Item(132, 13).Value = -34889438
$ws.Cells.Item(132, 14).Value = -16490.9999
$ws.Cells.Item(132, 11).Value = 34891968
$ws.Cells.Item(135, 8).Value = 792.88464
$ws.Cells.Item(135, 13).Value = -4063.8945
$ws.Cells.Item(135, 9).Value = 733.2105
$ws.Cells.Item(135, 11).Value = 6598.8945
$ws.Cells.Item(137, 11).Value = 596987.01
$ws.Cells.Item(137, 8).Value = 128489.93
$ws.Cells.Item(137, 12).Value = 4738.799999999999
$ws.Cells.Item(137, 9).Value = 198995.67
$ws.Cells.Item(137, 10).Value = 1579.6
$ws.Cells.Item(137, 13).Value = -594437.01
$ws.Cells.Item(137, 14).Value = -9838.799999999999
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(141, 11).Value = 31872.816
$ws.Cells.Item(141, 8).Value = 9297.462
$ws.Cells.Item(141, 13).Value = -26692.816
$ws.Cells.Item(141, 9).Value = 10624.272

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 11).Value = 3327822.8
$ws.Cells.Item(2, 8).Value = 2357643.5
$ws.Cells.Item(2, 12).Value = 1493.4286
$ws.Cells.Item(2, 9).Value = 3327822.8
$ws.Cells.Item(2, 10).Value = 1493.4286
$ws.Cells.Item(2, 13).Value = -3327709.8
$ws.Cells.Item(2, 14).Value = -1719.4286
$ws.Cells.Item(32, 11).Value = 1511.0804
$ws.Cells.Item(32, 8).Value = 2456.7605
$ws.Cells.Item(32, 13).Value = -1224.0804
$ws.Cells.Item(32, 9).Value = 1511.0804
$ws.Cells.Item(45, 8).Value = 7993835.5
$ws.Cells.Item(45, 9).Value = 9591642
$ws.Cells.Item(45, 13).Value = -9591265
$ws.Cells.Item(45, 11).Value = 9591642
$ws.Cells.Item(46, 13).ClearContents()
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 14).ClearContents()
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(61, 8).Value = 10428.866
$ws.Cells.Item(61, 13).Value = -10605
$ws.Cells.Item(61, 9).Value = 10817
$ws.Cells.Item(61, 11).Value = 10817
$ws.Cells.Item(74, 8).Value = 40267.63
$ws.Cells.Item(74, 13).Value = -6882.3486
$ws.Cells.Item(74, 9).Value = 7756.3486
$ws.Cells.Item(74, 11).Value = 7756.3486
$ws.Cells.Item(77, 11).Value = 38781.743
$ws.Cells.Item(77, 8).Value = 40267.63
$ws.Cells.Item(77, 13).Value = -34413.743
$ws.Cells.Item(77, 9).Value = 7756.3486
$ws.Cells.Item(97, 8).Value = 68830930
$ws.Cells.Item(97, 9).Value = 86030090
$ws.Cells.Item(97, 13).Value = -86029594
$ws.Cells.Item(97, 11).Value = 86030090
$ws.Cells.Item(110, 12).Value = 1004
$ws.Cells.Item(110, 10).Value = 1004
$ws.Cells.Item(110, 14).Value = -5094
$ws.Cells.Item(110, 8).Value = 2778644.2
$ws.Cells.Item(112, 12).Value = 20665.334
$ws.Cells.Item(112, 8).Value = 20665.334
$ws.Cells.Item(112, 10).Value = 20665.334
$ws.Cells.Item(112, 14).Value = -23619.334
$ws.Cells.Item(116, 8).Value = 2357643.5
$ws.Cells.Item(116, 12).Value = 1493.4286
$ws.Cells.Item(116, 9).Value = 3327822.8
$ws.Cells.Item(116, 10).Value = 1493.4286
$ws.Cells.Item(116, 13).Value = -3325528.8
$ws.Cells.Item(116, 14).Value = -6081.4286
$ws.Cells.Item(116, 11).Value = 3327822.8
$ws.Cells.Item(122, 11).Value = 6017.8668
$ws.Cells.Item(122, 8).Value = 338301.8
$ws.Cells.Item(122, 13).Value = -3567.8668
$ws.Cells.Item(122, 12).Value = 3685490.4
$ws.Cells.Item(122, 9).Value = 2005.9556
$ws.Cells.Item(122, 10).Value = 1228496.8
$ws.Cells.Item(122, 14).Value = -3690390.4
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 5109.5366
$ws.Cells.Item(132, 12).Value = 11998.2
$ws.Cells.Item(132, 9).Value = 5750
$ws.Cells.Item(132, 10).Value = 3999.4
$ws.Cells.Item(132, 13).Value = -14720
$ws.Cells.Item(132, 14).Value = -17058.2
$ws.Cells.Item(132, 11).Value = 17250
$ws.Cells.Item(136, 13).Value = -29901
$ws.Cells.Item(136, 9).Value = 10817
$ws.Cells.Item(136, 11).Value = 32451
$ws.Cells.Item(136, 8).Value = 10428.866
$ws.Cells.Item(138, 12).Value = 60000
$ws.Cells.Item(138, 8).Value = 60000
$ws.Cells.Item(138, 10).Value = 60000
$ws.Cells.Item(138, 14).Value = -70280
$ws.Cells.Item(139, 12).Value = 453331.88
$ws.Cells.Item(139, 10).Value = 453331.88
$ws.Cells.Item(139, 14).Value = -463611.88
$ws.Cells.Item(139, 8).Value = 458063.7

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 11).Value = 3327822.8
$ws.Cells.Item(3, 12).Value = 1493.4286
$ws.Cells.Item(3, 8).Value = 2357643.5
$ws.Cells.Item(3, 13).Value = -3327708.8
$ws.Cells.Item(3, 9).Value = 3327822.8
$ws.Cells.Item(3, 10).Value = 1493.4286
$ws.Cells.Item(3, 14).Value = -1721.4286
$ws.Cells.Item(22, 11).Value = 716.4
$ws.Cells.Item(22, 8).Value = 833.0909
$ws.Cells.Item(22, 13).Value = -543.4
$ws.Cells.Item(22, 9).Value = 716.4
$ws.Cells.Item(86, 11).Value = 7710448.5
$ws.Cells.Item(86, 8).Value = 5569701.5
$ws.Cells.Item(86, 12).Value = 3759.8
$ws.Cells.Item(86, 13).Value = -7709325.5
$ws.Cells.Item(86, 9).Value = 7710448.5
$ws.Cells.Item(86, 10).Value = 3759.8
$ws.Cells.Item(86, 14).Value = -6005.8
$ws.Cells.Item(89, 11).Value = 38552242.5
$ws.Cells.Item(89, 8).Value = 5569701.5
$ws.Cells.Item(89, 12).Value = 18799
$ws.Cells.Item(89, 9).Value = 7710448.5
$ws.Cells.Item(89, 10).Value = 3759.8
$ws.Cells.Item(89, 13).Value = -38546626.5
$ws.Cells.Item(89, 14).Value = -30031
$ws.Cells.Item(99, 11).Value = 15876402
$ws.Cells.Item(99, 8).Value = 8932187
$ws.Cells.Item(99, 12).Value = 3911.2856
$ws.Cells.Item(99, 9).Value = 15876402
$ws.Cells.Item(99, 10).Value = 3911.2856
$ws.Cells.Item(99, 13).Value = -15874904
$ws.Cells.Item(99, 14).Value = -6907.2856
$ws.Cells.Item(134, 11).Value = 20217.75
$ws.Cells.Item(134, 8).Value = 8071.16
$ws.Cells.Item(134, 9).Value = 6739.25
$ws.Cells.Item(134, 13).Value = -17682.75
$ws.Cells.Item(135, 8).Value = 72500
$ws.Cells.Item(135, 12).Value = 72500
$ws.Cells.Item(135, 10).Value = 72500
$ws.Cells.Item(135, 14).Value = -82640

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 12).Value = 1725
$ws.Cells.Item(4, 9).Value = 979.6667
$ws.Cells.Item(4, 10).Value = 1725
$ws.Cells.Item(4, 13).Value = -867.6667
$ws.Cells.Item(4, 14).Value = -1949
$ws.Cells.Item(4, 11).Value = 979.6667
$ws.Cells.Item(4, 8).Value = 1136.579
$ws.Cells.Item(12, 11).Value = 1096.1111
$ws.Cells.Item(12, 8).Value = 1179.1818
$ws.Cells.Item(12, 12).Value = 1553
$ws.Cells.Item(12, 9).Value = 1096.1111
$ws.Cells.Item(12, 10).Value = 1553
$ws.Cells.Item(12, 13).Value = -926.1111000000001
$ws.Cells.Item(12, 14).Value = -1893
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 14).ClearContents()
$ws.Cells.Item(17, 8).Value = 200000
$ws.Cells.Item(22, 11).Value = 1099.6
$ws.Cells.Item(22, 12).Value = 2887.5557
$ws.Cells.Item(22, 8).Value = 2249
$ws.Cells.Item(22, 13).Value = -749.5999999999999
$ws.Cells.Item(22, 9).Value = 1099.6
$ws.Cells.Item(22, 10).Value = 2887.5557
$ws.Cells.Item(22, 14).Value = -3587.5557
$ws.Cells.Item(31, 11).Value = 953.7241
$ws.Cells.Item(31, 13).Value = -658.7241
$ws.Cells.Item(31, 8).Value = 3269.04
$ws.Cells.Item(31, 12).Value = 4214.7324
$ws.Cells.Item(31, 9).Value = 953.7241
$ws.Cells.Item(31, 10).Value = 4214.7324
$ws.Cells.Item(31, 14).Value = -4804.7324
$ws.Cells.Item(34, 14).Value = -4618.7324
$ws.Cells.Item(34, 11).Value = 953.7241
$ws.Cells.Item(34, 8).Value = 3269.04
$ws.Cells.Item(34, 12).Value = 4214.7324
$ws.Cells.Item(34, 13).Value = -751.7241
$ws.Cells.Item(34, 9).Value = 953.7241
$ws.Cells.Item(34, 10).Value = 4214.7324
$ws.Cells.Item(58, 8).Value = 3423.625
$ws.Cells.Item(58, 13).Value = -2962.6667
$ws.Cells.Item(58, 9).Value = 3165.6667
$ws.Cells.Item(58, 11).Value = 3165.6667
$ws.Cells.Item(99, 11).Value = 3509.0557
$ws.Cells.Item(99, 8).Value = 3652.3794
$ws.Cells.Item(99, 12).Value = 3886.9092
$ws.Cells.Item(99, 9).Value = 3509.0557
$ws.Cells.Item(99, 10).Value = 3886.9092
$ws.Cells.Item(99, 13).Value = -2011.0557
$ws.Cells.Item(99, 14).Value = -6882.9092
$ws.Cells.Item(107, 9).Value = 797
$ws.Cells.Item(107, 13).Value = 1123
$ws.Cells.Item(107, 11).Value = 797
$ws.Cells.Item(107, 8).Value = 947.375
$ws.Cells.Item(122, 11).Value = 11621.25
$ws.Cells.Item(122, 8).Value = 3969.875
$ws.Cells.Item(122, 13).Value = -9171.25
$ws.Cells.Item(122, 12).Value = 12005.7501
$ws.Cells.Item(122, 9).Value = 3873.75
$ws.Cells.Item(122, 10).Value = 4001.9167
$ws.Cells.Item(122, 14).Value = -16905.7501
$ws.Cells.Item(126, 8).Value = 3652.3794
$ws.Cells.Item(126, 12).Value = 11660.7276
$ws.Cells.Item(126, 9).Value = 3509.0557
$ws.Cells.Item(126, 10).Value = 3886.9092
$ws.Cells.Item(126, 13).Value = -8057.167099999999
$ws.Cells.Item(126, 14).Value = -16600.7276
$ws.Cells.Item(126, 11).Value = 10527.1671
$ws.Cells.Item(132, 8).Value = 43948.832
$ws.Cells.Item(132, 9).Value = 52067.9
$ws.Cells.Item(132, 13).Value = -153673.7
$ws.Cells.Item(132, 11).Value = 156203.7
$ws.Cells.Item(134, 11).Value = 14691.483
$ws.Cells.Item(134, 8).Value = 7167.8335
$ws.Cells.Item(134, 12).Value = 40701
$ws.Cells.Item(134, 9).Value = 4897.161
$ws.Cells.Item(134, 10).Value = 13567
$ws.Cells.Item(134, 13).Value = -12156.483
$ws.Cells.Item(134, 14).Value = -45771
$ws.Cells.Item(136, 13).Value = -6947.000100000001
$ws.Cells.Item(136, 9).Value = 3165.6667
$ws.Cells.Item(136, 11).Value = 9497.000100000001
$ws.Cells.Item(136, 8).Value = 3423.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 14).Value = -102332.996
$ws.Cells.Item(5, 12).Value = 102108.996
$ws.Cells.Item(5, 11).Value = 1613.0625
$ws.Cells.Item(5, 8).Value = 5826.9473
$ws.Cells.Item(5, 13).Value = -1501.0625
$ws.Cells.Item(5, 9).Value = 537.6875
$ws.Cells.Item(5, 10).Value = 34036.332
$ws.Cells.Item(23, 8).Value = 201.61539
$ws.Cells.Item(23, 12).Value = 746.6666700000001
$ws.Cells.Item(23, 9).Value = 95.25
$ws.Cells.Item(23, 10).Value = 248.88889
$ws.Cells.Item(23, 13).Value = -50.75
$ws.Cells.Item(23, 14).Value = -1216.66667
$ws.Cells.Item(23, 11).Value = 285.75
$ws.Cells.Item(39, 12).Value = 14970
$ws.Cells.Item(39, 10).Value = 4990
$ws.Cells.Item(39, 14).Value = -15558
$ws.Cells.Item(39, 8).Value = 3579.6667
$ws.Cells.Item(47, 11).Value = 1323
$ws.Cells.Item(47, 8).Value = 510.875
$ws.Cells.Item(47, 12).Value = 3000
$ws.Cells.Item(47, 9).Value = 441
$ws.Cells.Item(47, 10).Value = 1000
$ws.Cells.Item(47, 13).Value = -892
$ws.Cells.Item(47, 14).Value = -3862
$ws.Cells.Item(55, 8).Value = 61190.59
$ws.Cells.Item(55, 13).Value = -1698
$ws.Cells.Item(55, 12).Value = 443138.58
$ws.Cells.Item(55, 9).Value = 625
$ws.Cells.Item(55, 10).Value = 147712.86
$ws.Cells.Item(55, 14).Value = -443492.58
$ws.Cells.Item(55, 11).Value = 1875
$ws.Cells.Item(56, 11).Value = 7121.3477
$ws.Cells.Item(56, 8).Value = 7121.3477
$ws.Cells.Item(56, 13).Value = -6591.3477
$ws.Cells.Item(56, 9).Value = 7121.3477
$ws.Cells.Item(60, 11).Value = 5955.9231
$ws.Cells.Item(60, 8).Value = 2199.9375
$ws.Cells.Item(60, 12).Value = 9390
$ws.Cells.Item(60, 13).Value = -5704.9231
$ws.Cells.Item(60, 9).Value = 1985.3077
$ws.Cells.Item(60, 10).Value = 3130
$ws.Cells.Item(60, 14).Value = -9892
$ws.Cells.Item(68, 8).Value = 1680.862
$ws.Cells.Item(68, 12).Value = 5602.5
$ws.Cells.Item(68, 9).Value = 1506.6666
$ws.Cells.Item(68, 10).Value = 1867.5
$ws.Cells.Item(68, 13).Value = -3708.9998
$ws.Cells.Item(68, 14).Value = -7224.5
$ws.Cells.Item(68, 11).Value = 4519.9998
$ws.Cells.Item(71, 8).Value = 1680.862
$ws.Cells.Item(71, 13).Value = -9503.999400000001
$ws.Cells.Item(71, 12).Value = 16807.5
$ws.Cells.Item(71, 9).Value = 1506.6666
$ws.Cells.Item(71, 10).Value = 1867.5
$ws.Cells.Item(71, 14).Value = -24919.5
$ws.Cells.Item(71, 11).Value = 13559.9994
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 8).Value = 1666.6666
$ws.Cells.Item(86, 12).Value = 4999.9998
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 1666.6666
$ws.Cells.Item(86, 14).Value = -7371.9998
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 8).Value = 1666.6666
$ws.Cells.Item(89, 12).Value = 14999.9994
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 1666.6666
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(89, 14).Value = -26855.9994
$ws.Cells.Item(107, 12).Value = 3000
$ws.Cells.Item(107, 9).Value = 550.125
$ws.Cells.Item(107, 10).Value = 1000
$ws.Cells.Item(107, 13).Value = 269.625
$ws.Cells.Item(107, 14).Value = -6840
$ws.Cells.Item(107, 11).Value = 1650.375
$ws.Cells.Item(107, 8).Value = 640.1
$ws.Cells.Item(113, 8).Value = 2544.9678
$ws.Cells.Item(113, 12).Value = 5980.0908
$ws.Cells.Item(113, 9).Value = 3893.3333
$ws.Cells.Item(113, 10).Value = 1993.3636
$ws.Cells.Item(113, 13).Value = -9509.999899999999
$ws.Cells.Item(113, 14).Value = -10320.0908
$ws.Cells.Item(113, 11).Value = 11679.9999
$ws.Cells.Item(131, 11).Value = 5280.6
$ws.Cells.Item(131, 8).Value = 5365.8945
$ws.Cells.Item(131, 12).Value = 28116.669
$ws.Cells.Item(131, 13).Value = -240.6000000000004
$ws.Cells.Item(131, 9).Value = 1760.2
$ws.Cells.Item(131, 10).Value = 9372.223
$ws.Cells.Item(131, 14).Value = -38196.669
$ws.Cells.Item(132, 8).Value = 2673.75
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 9).Value = 2673.75
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 13).Value = -21533.75
$ws.Cells.Item(132, 14).ClearContents()
$ws.Cells.Item(132, 11).Value = 24063.75
$ws.Cells.Item(133, 12).Value = 33999.999
$ws.Cells.Item(133, 10).Value = 11333.333
$ws.Cells.Item(133, 14).Value = -44119.999
$ws.Cells.Item(133, 8).Value = 6885.6665
$ws.Cells.Item(134, 11).Value = 4645.071599999999
$ws.Cells.Item(134, 8).Value = 2445.1333
$ws.Cells.Item(134, 12).Value = 45000
$ws.Cells.Item(134, 9).Value = 1548.3572
$ws.Cells.Item(134, 10).Value = 15000
$ws.Cells.Item(134, 13).Value = 424.9284000000007
$ws.Cells.Item(134, 14).Value = -55140
$ws.Cells.Item(135, 8).Value = 5826.9473
$ws.Cells.Item(135, 13).Value = -2304.1875
$ws.Cells.Item(135, 12).Value = 306326.988
$ws.Cells.Item(135, 9).Value = 537.6875
$ws.Cells.Item(135, 10).Value = 34036.332
$ws.Cells.Item(135, 14).Value = -311396.988
$ws.Cells.Item(135, 11).Value = 4839.1875
$ws.Cells.Item(136, 13).Value = -1889.700000000001
$ws.Cells.Item(136, 9).Value = 2329.9
$ws.Cells.Item(136, 11).Value = 6989.700000000001
$ws.Cells.Item(136, 8).Value = 2329.9
$ws.Cells.Item(137, 11).Value = 9338.571599999999
$ws.Cells.Item(137, 8).Value = 4175.067
$ws.Cells.Item(137, 12).Value = 15313.5
$ws.Cells.Item(137, 9).Value = 3112.8572
$ws.Cells.Item(137, 10).Value = 5104.5
$ws.Cells.Item(137, 13).Value = -4238.571599999999
$ws.Cells.Item(137, 14).Value = -25513.5
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 9).Value = 1370.15
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 13).Value = 1029.549999999999
$ws.Cells.Item(139, 14).ClearContents()
$ws.Cells.Item(139, 11).Value = 4110.450000000001
$ws.Cells.Item(139, 8).Value = 1370.15
$ws.Cells.Item(140, 11).Value = 3407.8236
$ws.Cells.Item(140, 8).Value = 1817.8695
$ws.Cells.Item(140, 13).Value = 1772.1764
$ws.Cells.Item(140, 9).Value = 1135.9412

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 14).Value = -20346
$ws.Cells.Item(24, 11).Value = 17999.8
$ws.Cells.Item(24, 12).Value = 20000
$ws.Cells.Item(24, 8).Value = 18333.166
$ws.Cells.Item(24, 13).Value = -17826.8
$ws.Cells.Item(24, 9).Value = 17999.8
$ws.Cells.Item(24, 10).Value = 20000
$ws.Cells.Item(49, 12).Value = 34750
$ws.Cells.Item(49, 10).Value = 34750
$ws.Cells.Item(49, 14).Value = -35118
$ws.Cells.Item(49, 8).Value = 34750
$ws.Cells.Item(80, 8).Value = 43884576
$ws.Cells.Item(80, 13).Value = -131122362
$ws.Cells.Item(80, 12).Value = 265186.12
$ws.Cells.Item(80, 9).Value = 131123360
$ws.Cells.Item(80, 10).Value = 265186.12
$ws.Cells.Item(80, 14).Value = -267182.12
$ws.Cells.Item(80, 11).Value = 131123360
$ws.Cells.Item(83, 11).Value = 655616800
$ws.Cells.Item(83, 8).Value = 43884576
$ws.Cells.Item(83, 12).Value = 1325930.6
$ws.Cells.Item(83, 9).Value = 131123360
$ws.Cells.Item(83, 10).Value = 265186.12
$ws.Cells.Item(83, 13).Value = -655611808
$ws.Cells.Item(83, 14).Value = -1335914.6
$ws.Cells.Item(103, 8).Value = 44868
$ws.Cells.Item(103, 12).Value = 44868
$ws.Cells.Item(103, 10).Value = 44868
$ws.Cells.Item(103, 14).Value = -47212
$ws.Cells.Item(126, 8).Value = 2929834.8
$ws.Cells.Item(126, 9).Value = 1468747
$ws.Cells.Item(126, 13).Value = -4403771
$ws.Cells.Item(126, 11).Value = 4406241
$ws.Cells.Item(132, 8).Value = 5386.45
$ws.Cells.Item(132, 12).Value = 50076
$ws.Cells.Item(132, 9).Value = 2988.303
$ws.Cells.Item(132, 10).Value = 16692
$ws.Cells.Item(132, 13).Value = -6434.909
$ws.Cells.Item(132, 14).Value = -55136
$ws.Cells.Item(132, 11).Value = 8964.909
$ws.Cells.Item(139, 12).Value = 51000
$ws.Cells.Item(139, 10).Value = 51000
$ws.Cells.Item(139, 14).Value = -61280
$ws.Cells.Item(139, 8).Value = 51000

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 11).Value = 15002500
$ws.Cells.Item(2, 8).Value = 12005000
$ws.Cells.Item(2, 9).Value = 15002500
$ws.Cells.Item(2, 13).Value = -15002388
$ws.Cells.Item(9, 11).Value = 375
$ws.Cells.Item(9, 8).Value = 1299.1666
$ws.Cells.Item(9, 9).Value = 375
$ws.Cells.Item(9, 13).Value = -151
$ws.Cells.Item(19, 11).Value = 650
$ws.Cells.Item(19, 8).Value = 480
$ws.Cells.Item(19, 13).Value = -480
$ws.Cells.Item(19, 12).Value = 140
$ws.Cells.Item(19, 10).Value = 140
$ws.Cells.Item(19, 9).Value = 650
$ws.Cells.Item(19, 14).Value = -480
$ws.Cells.Item(22, 11).Value = 179036.4
$ws.Cells.Item(22, 12).Value = 1574.75
$ws.Cells.Item(22, 8).Value = 100164.555
$ws.Cells.Item(22, 13).Value = -178741.4
$ws.Cells.Item(22, 9).Value = 179036.4
$ws.Cells.Item(22, 10).Value = 1574.75
$ws.Cells.Item(22, 14).Value = -2164.75
$ws.Cells.Item(27, 13).Value = -178929.4
$ws.Cells.Item(27, 9).Value = 179036.4
$ws.Cells.Item(27, 10).Value = 1574.75
$ws.Cells.Item(27, 14).Value = -1788.75
$ws.Cells.Item(27, 11).Value = 179036.4
$ws.Cells.Item(27, 8).Value = 100164.555
$ws.Cells.Item(27, 12).Value = 1574.75
$ws.Cells.Item(40, 14).Value = -10280.8
$ws.Cells.Item(40, 11).Value = 7154.8125
$ws.Cells.Item(40, 12).Value = 10008.8
$ws.Cells.Item(40, 8).Value = 8252.5
$ws.Cells.Item(40, 13).Value = -7018.8125
$ws.Cells.Item(40, 9).Value = 7154.8125
$ws.Cells.Item(40, 10).Value = 10008.8
$ws.Cells.Item(42, 8).Value = 8612.5
$ws.Cells.Item(42, 9).Value = 8612.5
$ws.Cells.Item(42, 13).Value = -8049.5
$ws.Cells.Item(42, 11).Value = 8612.5
$ws.Cells.Item(46, 13).Value = -1479.3334
$ws.Cells.Item(46, 9).Value = 1667.3334
$ws.Cells.Item(46, 10).Value = 1899.6666
$ws.Cells.Item(46, 14).Value = -2275.6666
$ws.Cells.Item(46, 11).Value = 1667.3334
$ws.Cells.Item(46, 8).Value = 1783.5
$ws.Cells.Item(46, 12).Value = 1899.6666
$ws.Cells.Item(49, 9).Value = 8612.5
$ws.Cells.Item(49, 13).Value = -8465.5
$ws.Cells.Item(49, 11).Value = 8612.5
$ws.Cells.Item(49, 8).Value = 8612.5
$ws.Cells.Item(55, 8).Value = 1985.762
$ws.Cells.Item(55, 13).Value = -1179.3636
$ws.Cells.Item(55, 12).Value = 2682.5
$ws.Cells.Item(55, 9).Value = 1352.3636
$ws.Cells.Item(55, 10).Value = 2682.5
$ws.Cells.Item(55, 14).Value = -3028.5
$ws.Cells.Item(55, 11).Value = 1352.3636
$ws.Cells.Item(61, 8).Value = 6946142
$ws.Cells.Item(61, 13).Value = -9260700
$ws.Cells.Item(61, 9).Value = 9260902
$ws.Cells.Item(61, 11).Value = 9260902
$ws.Cells.Item(82, 8).Value = 65973612
$ws.Cells.Item(82, 12).Value = 1057
$ws.Cells.Item(82, 10).Value = 1057
$ws.Cells.Item(82, 14).Value = -1779
$ws.Cells.Item(85, 14).Value = -3553
$ws.Cells.Item(85, 12).Value = 1057
$ws.Cells.Item(85, 8).Value = 65973612
$ws.Cells.Item(85, 10).Value = 1057
$ws.Cells.Item(110, 12).Value = 64833.332
$ws.Cells.Item(110, 10).Value = 64833.332
$ws.Cells.Item(110, 14).Value = -73013.33199999999
$ws.Cells.Item(110, 8).Value = 64833.332
$ws.Cells.Item(113, 8).Value = 6946142
$ws.Cells.Item(113, 9).Value = 9260902
$ws.Cells.Item(113, 13).Value = -9258732
$ws.Cells.Item(113, 11).Value = 9260902
$ws.Cells.Item(122, 11).Value = 12114.4284
$ws.Cells.Item(122, 8).Value = 5603.5654
$ws.Cells.Item(122, 13).Value = -9664.428400000001
$ws.Cells.Item(122, 9).Value = 4038.1428
$ws.Cells.Item(132, 8).Value = 3080.4119
$ws.Cells.Item(132, 12).Value = 12620.25
$ws.Cells.Item(132, 9).Value = 2733.8462
$ws.Cells.Item(132, 10).Value = 4206.75
$ws.Cells.Item(132, 13).Value = -5671.5386
$ws.Cells.Item(132, 14).Value = -17680.25
$ws.Cells.Item(132, 11).Value = 8201.5386
$ws.Cells.Item(136, 13).Value = -101764.8
$ws.Cells.Item(136, 9).Value = 34771.6
$ws.Cells.Item(136, 10).Value = 4269.7
$ws.Cells.Item(136, 14).Value = -17909.1
$ws.Cells.Item(136, 11).Value = 104314.8
$ws.Cells.Item(136, 8).Value = 22570.84
$ws.Cells.Item(136, 12).Value = 12809.1
$ws.Cells.Item(137, 8).Value = 50429
$ws.Cells.Item(137, 12).Value = 50429
$ws.Cells.Item(137, 10).Value = 50429
$ws.Cells.Item(137, 14).Value = -60629

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 11).Value = 7008.8334
$ws.Cells.Item(122, 8).Value = 2222.276
$ws.Cells.Item(122, 13).Value = -4558.8334
$ws.Cells.Item(122, 12).Value = 6107.1819
$ws.Cells.Item(122, 9).Value = 2336.2778
$ws.Cells.Item(122, 10).Value = 2035.7273
$ws.Cells.Item(122, 14).Value = -11007.1819
$ws.Cells.Item(132, 8).Value = 21066928
$ws.Cells.Item(132, 12).Value = 3303927.6
$ws.Cells.Item(132, 9).Value = 26321038
$ws.Cells.Item(132, 10).Value = 1101309.2
$ws.Cells.Item(132, 13).Value = -78960584
$ws.Cells.Item(132, 14).Value = -3308987.6
$ws.Cells.Item(132, 11).Value = 78963114
